$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.320918321609497
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 3.240891695022583
$ws.Range("D1").Value = 1.611703395843506
$ws.Range("E1").Value = 1.188177466392517
